$wb = $excel.ActiveWorkbook

# --- Sheet "pet_gourmet" ---
$ws1 = $wb.Worksheets.Item("pet_gourmet")

$ws1.Range("E11").Value = 810
$ws1.Range("F11").Value = 10

$ws1.Range("C12").Value = 10
$ws1.Range("E12").Value = 3000
$ws1.Range("F12").Value = 1510

$ws1.Range("C13").Value = 1510
$ws1.Range("E13").Value = 0

$ws1.Range("E14").Value = 1810
$ws1.Range("F14").Value = 1010

$ws1.Range("C15").Value = 1010
$ws1.Range("E15").Value = 0
$ws1.Range("F15").Value = 10

$ws1.Range("C16").Value = 10
$ws1.Range("E16").Value = 1000

# --- Sheet "patas_pack" ---
$ws2 = $wb.Worksheets.Item("patas_pack")

$ws2.Range("D11").Value = 810
$ws2.Range("E11").Value = 710

$ws2.Range("D12").Value = 3000
$ws2.Range("E12").Value = 3000

$ws2.Range("D13").Value = 0
$ws2.Range("E13").Value = 0

$ws2.Range("D14").Value = 1810
$ws2.Range("E14").Value = 1710

$ws2.Range("D15").Value = 0
$ws2.Range("E15").Value = 0

$ws2.Range("D16").Value = 1000
$ws2.Range("E16").Value = 1000
